$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Tue Sep 26 21:30:29 EDT 2023"
$ws.Range("B3").Value = "Tue Sep 26 21:30:42 EDT 2023"
$ws.Range("B4").Value = "Tue Sep 26 21:30:56 EDT 2023"
